# Apply the epexspot_prices.xlsx update:
#  1. "Prix Spot": insert a new date column "07-dec" right before the
#     existing "01-oct." column (which shifts EJ:FN -> EK:FO), and fill the
#     new column's hourly data rows (2-25) with "-" placeholders.
#  2. "Gaz": append a new row with the 2025-12-05 Last Price.
#  3. "CO2": append a new row with the 2025-12-05 Last Price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Prix Spot" sheet - insert the "07-dec" column before column EJ
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Columns("EJ").Insert()

$ws.Range("EJ1").Value = "07-dec"
$ws.Range("EJ2:EJ25").Value = "-"

# ---------------------------------------------------------------------
# 2) "Gaz" sheet - append row 170
# ---------------------------------------------------------------------
$gaz = $wb.Worksheets.Item("Gaz")

$gaz.Range("A170").NumberFormat = "@"
$gaz.Range("A170").Value = "2025-12-05"
$gaz.Range("A170").ClearFormats()
$gaz.Range("B170").Value = 25.965

# ---------------------------------------------------------------------
# 3) "CO2" sheet - append row 170
# ---------------------------------------------------------------------
$co2 = $wb.Worksheets.Item("CO2")

$co2.Range("A170").NumberFormat = "@"
$co2.Range("A170").Value = "2025-12-05"
$co2.Range("A170").ClearFormats()
$co2.Range("B170").Value = 81.78
